# fix bug in parameters
#
# The "UK" policy schedule only listed a row for each EUROMOD parameter
# file that happened to already have a Policy_Start_Year / Policy_
# System_Year filled in (uk_2015_std.txt and uk_2019_std.txt). That
# was a bug: every available uk_*_std.txt parameter file (2011-2026)
# needs its own row, even though most of them don't have an explicit
# start/system year recorded yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")

# Years that already carry a Policy_Start_Year / Policy_System_Year
# value; everything else just gets the filename (blank years, same as
# the rest of the table).
$yearsWithValues = @{ 2015 = "2015"; 2019 = "2019" }

$firstYear = 2011
$lastYear = 2026
$startRow = 2
$endRow = $startRow + ($lastYear - $firstYear)

# Start from a clean slate so stale values from the old (shorter)
# table don't leak into rows that should now be blank.
$dataRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 4))
$dataRange.ClearContents()

# Make sure the Policy_Start_Year / Policy_System_Year columns are
# written as text (matching the existing "2015"/"2019" text entries)
# rather than being auto-converted to numbers.
$yearColumns = $ws.Range($ws.Cells.Item($startRow, 2), $ws.Cells.Item($endRow, 3))
$yearColumns.NumberFormat = "@"

$row = $startRow
for ($year = $firstYear; $year -le $lastYear; $year++) {
    $ws.Cells.Item($row, 1).Value = "uk_" + $year + "_std.txt"

    if ($yearsWithValues.ContainsKey($year)) {
        $ws.Cells.Item($row, 2).Value = $yearsWithValues[$year]
        $ws.Cells.Item($row, 3).Value = $yearsWithValues[$year]
    }

    $row = $row + 1
}

# Drop the temporary text formatting again so the cells fall back to
# the default (unstyled) look used by the rest of the sheet.
$yearColumns.ClearFormats()
